$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5393.5557
$ws.Range("I100").Value = 2950.6924
$ws.Range("K100").Value = 2950.6924
$ws.Range("M100").Value = -2409.6924
$ws.Range("H129").Value = 2021.3077
$ws.Range("J129").Value = 2333
$ws.Range("L129").Value = 6999
$ws.Range("N129").Value = -16999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2554.3333
$ws.Range("I102").Value = 2182.05
$ws.Range("K102").Value = 2182.05
$ws.Range("M102").Value = -560.0500000000002
$ws.Range("H122").Value = 9261533
$ws.Range("I122").Value = 18519950
$ws.Range("K122").Value = 55559850
$ws.Range("M122").Value = -55557400
$ws.Range("H132").Value = 3421.25
$ws.Range("I132").Value = 2905.5
$ws.Range("K132").Value = 8716.5
$ws.Range("M132").Value = -6186.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1346.8
$ws.Range("I11").Value = 433.5
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 433.5
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = -293.5
$ws.Range("N11").Value = -5280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2479.8
$ws.Range("I99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("H105").Value = 2461.5
$ws.Range("I105").Value = 2461.5
$ws.Range("K105").Value = 2461.5
$ws.Range("M105").Value = -714.5
$ws.Range("H126").Value = 2479.8
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 2425.96
$ws.Range("I132").Value = 2425.96
$ws.Range("K132").Value = 7277.88
$ws.Range("M132").Value = -4747.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 9182035
$ws.Range("I17").Value = 11111348
$ws.Range("K17").Value = 33334044
$ws.Range("M17").Value = -33333875
$ws.Range("H34").Value = 3674.611
$ws.Range("J34").Value = 3972.4546
$ws.Range("L34").Value = 11917.3638
$ws.Range("N34").Value = -12085.3638
$ws.Range("H39").Value = 10004
$ws.Range("J39").Value = 10004
$ws.Range("L39").Value = 30012
$ws.Range("N39").Value = -30600
$ws.Range("H63").Value = 938.5
$ws.Range("I63").Value = 938.5
$ws.Range("K63").Value = 2815.5
$ws.Range("M63").Value = -2066.5
$ws.Range("H66").Value = 938.5
$ws.Range("I66").Value = 938.5
$ws.Range("K66").Value = 8446.5
$ws.Range("M66").Value = -4702.5
$ws.Range("H70").Value = 5148
$ws.Range("I70").Value = 4722
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 14166
$ws.Range("L70").Value = 18000
$ws.Range("M70").Value = -13851
$ws.Range("N70").Value = -18630
$ws.Range("H73").Value = 5148
$ws.Range("I73").Value = 4722
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 14166
$ws.Range("L73").Value = 18000
$ws.Range("M73").Value = -13074
$ws.Range("N73").Value = -20184
$ws.Range("H113").Value = 3562.25
$ws.Range("I113").Value = 3375
$ws.Range("J113").Value = 3749.5
$ws.Range("K113").Value = 10125
$ws.Range("L113").Value = 11248.5
$ws.Range("M113").Value = -7955
$ws.Range("N113").Value = -15588.5
$ws.Range("H122").Value = 2541
$ws.Range("H129").Value = 2827.3
$ws.Range("J129").Value = 3720.7856
$ws.Range("L129").Value = 11162.3568
$ws.Range("N129").Value = -21162.3568
$ws.Range("H132").Value = 2391
$ws.Range("I132").Value = 1579
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 14211
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -11681
$ws.Range("N132").Value = -32060
$ws.Range("H134").Value = 3157.68
$ws.Range("I134").Value = 1262.8
$ws.Range("K134").Value = 3788.4
$ws.Range("M134").Value = 1281.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1149.6666
$ws.Range("I6").Value = 1524.5
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 1524.5
$ws.Range("L6").Value = 400
$ws.Range("M6").Value = -1411.5
$ws.Range("N6").Value = -626
$ws.Range("H11").Value = 12501377
$ws.Range("I11").Value = 13540242
$ws.Range("J11").Value = 35000
$ws.Range("K11").Value = 13540242
$ws.Range("L11").Value = 35000
$ws.Range("M11").Value = -13540103
$ws.Range("N11").Value = -35278
$ws.Range("H16").Value = 1149.6666
$ws.Range("I16").Value = 1524.5
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 1524.5
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -1274.5
$ws.Range("N16").Value = -900
$ws.Range("H17").Value = 403
$ws.Range("I17").Value = 299
$ws.Range("J17").Value = 507
$ws.Range("K17").Value = 299
$ws.Range("L17").Value = 507
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -843
$ws.Range("H63").Value = 28499
$ws.Range("J63").Value = 28499
$ws.Range("L63").Value = 28499
$ws.Range("N63").Value = -29871
$ws.Range("H66").Value = 28499
$ws.Range("J66").Value = 28499
$ws.Range("L66").Value = 85497
$ws.Range("N66").Value = -92361
$ws.Range("H122").Value = 10235.518
$ws.Range("I122").Value = 10418.84
$ws.Range("J122").Value = 9089.75
$ws.Range("K122").Value = 31256.52
$ws.Range("L122").Value = 27269.25
$ws.Range("M122").Value = -28806.52
$ws.Range("N122").Value = -32169.25
$ws.Range("H123").Value = 27142.285
$ws.Range("J123").Value = 27142.285
$ws.Range("L123").Value = 27142.285
$ws.Range("N123").Value = -32042.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6013
$ws.Range("I7").Value = 5633.1665
$ws.Range("K7").Value = 5633.1665
$ws.Range("M7").Value = -5521.1665
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 8042.909
$ws.Range("I40").Value = 6337.091
$ws.Range("K40").Value = 6337.091
$ws.Range("M40").Value = -6201.091
$ws.Range("H46").Value = 1736.7273
$ws.Range("J46").Value = 1938
$ws.Range("L46").Value = 1938
$ws.Range("N46").Value = -2314
$ws.Range("H126").Value = 6013
$ws.Range("I126").Value = 5633.1665
$ws.Range("K126").Value = 16899.4995
$ws.Range("M126").Value = -14429.4995
$ws.Range("H132").Value = 5845.1577
$ws.Range("I132").Value = 5857.923
$ws.Range("J132").Value = 5817.5
$ws.Range("K132").Value = 17573.769
$ws.Range("L132").Value = 17452.5
$ws.Range("M132").Value = -15043.769
$ws.Range("N132").Value = -22512.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 185000
$ws.Range("J64").Value = 185000
$ws.Range("L64").Value = 185000
$ws.Range("N64").Value = -185496
$ws.Range("H67").Value = 185000
$ws.Range("J67").Value = 185000
$ws.Range("L67").Value = 185000
$ws.Range("N67").Value = -186716
$ws.Range("H113").Value = 2286.625
$ws.Range("I113").Value = 2074.5
$ws.Range("K113").Value = 6223.5
$ws.Range("M113").Value = -4053.5
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
